$wb = $excel.ActiveWorkbook

# "Generate Report for Handback": refresh the "Latest Handback DateTime" (column K)
# on the zh-cn and de-de sheets for the 4fef9958-... row (row 2).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("K2").Value = "2016-10-24 10:04:43"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-10-24 10:04:59"
